$wb = $excel.ActiveWorkbook

# ALC row 5
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 52.6
$ws.Range("I5").Value = 48.444443
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 48.444443
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = 66.55555699999999
$ws.Range("N5").Value = -320

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1442.6666
$ws.Range("I40").Value = 1332.6666
$ws.Range("J40").Value = 1589.3334
$ws.Range("K40").Value = 1332.6666
$ws.Range("L40").Value = 1589.3334
$ws.Range("M40").Value = -1157.6666
$ws.Range("N40").Value = -1939.3334

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2371.4285
$ws.Range("I51").Value = 1399.5
$ws.Range("J51").Value = 2760.2
$ws.Range("K51").Value = 1399.5
$ws.Range("L51").Value = 2760.2
$ws.Range("M51").Value = -915.5
$ws.Range("N51").Value = -3728.2

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3122.7273
$ws.Range("I64").Value = 2743.75
$ws.Range("J64").Value = 4133.3335
$ws.Range("K64").Value = 2743.75
$ws.Range("L64").Value = 4133.3335
$ws.Range("M64").Value = -2495.75
$ws.Range("N64").Value = -4629.3335

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3122.7273
$ws.Range("I67").Value = 2743.75
$ws.Range("J67").Value = 4133.3335
$ws.Range("K67").Value = 2743.75
$ws.Range("L67").Value = 4133.3335
$ws.Range("M67").Value = -1885.75
$ws.Range("N67").Value = -5849.3335

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3518.0667
$ws.Range("I74").Value = 3325.5454
$ws.Range("J74").Value = 4047.5
$ws.Range("K74").Value = 3325.5454
$ws.Range("L74").Value = 4047.5
$ws.Range("M74").Value = -2389.5454
$ws.Range("N74").Value = -5919.5

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3518.0667
$ws.Range("I77").Value = 3325.5454
$ws.Range("J77").Value = 4047.5
$ws.Range("K77").Value = 16627.727
$ws.Range("L77").Value = 20237.5
$ws.Range("M77").Value = -11947.727
$ws.Range("N77").Value = -29597.5

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 10754459
$ws.Range("I100").Value = 12821690
$ws.Range("J100").Value = 4857.2
$ws.Range("K100").Value = 12821690
$ws.Range("L100").Value = 4857.2
$ws.Range("M100").Value = -12821149
$ws.Range("N100").Value = -5939.2

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2336.7273
$ws.Range("I106").Value = 2220.4
$ws.Range("J106").Value = 3500
$ws.Range("K106").Value = 2220.4
$ws.Range("L106").Value = 3500
$ws.Range("M106").Value = -1589.4
$ws.Range("N106").Value = -4762

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2036.8125
$ws.Range("J113").Value = 1839.8
$ws.Range("L113").Value = 1839.8
$ws.Range("N113").Value = -8347.799999999999

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 879.25
$ws.Range("J129").Value = 1370
$ws.Range("L129").Value = 4110
$ws.Range("N129").Value = -14110

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 26316840
$ws.Range("I132").Value = 32258744
$ws.Range("J132").Value = 2692.2856
$ws.Range("K132").Value = 96776232
$ws.Range("L132").Value = 8076.8568
$ws.Range("M132").Value = -96773702
$ws.Range("N132").Value = -13136.8568

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1094.6842
$ws.Range("I137").Value = 907.13336
$ws.Range("J137").Value = 1798
$ws.Range("K137").Value = 2721.40008
$ws.Range("L137").Value = 5394
$ws.Range("M137").Value = -171.4000800000003
$ws.Range("N137").Value = -10494

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1862.381
$ws.Range("I2").Value = 883.38464
$ws.Range("K2").Value = 883.38464
$ws.Range("M2").Value = -770.38464

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 355667.03
$ws.Range("I32").Value = 2529.554
$ws.Range("J32").Value = 4088834.5
$ws.Range("K32").Value = 2529.554
$ws.Range("L32").Value = 4088834.5
$ws.Range("M32").Value = -2242.554
$ws.Range("N32").Value = -4089408.5

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3406.6086
$ws.Range("I45").Value = 3395.5715
$ws.Range("K45").Value = 3395.5715
$ws.Range("M45").Value = -3018.5715

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3796.3157
$ws.Range("I102").Value = 2347.182
$ws.Range("J102").Value = 5788.875
$ws.Range("K102").Value = 2347.182
$ws.Range("L102").Value = 5788.875
$ws.Range("M102").Value = -725.1819999999998
$ws.Range("N102").Value = -9032.875

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 717.0323
$ws.Range("I110").Value = 713.26086
$ws.Range("J110").Value = 727.875
$ws.Range("K110").Value = 713.26086
$ws.Range("L110").Value = 727.875
$ws.Range("M110").Value = 1331.73914
$ws.Range("N110").Value = -4817.875

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1862.381
$ws.Range("I116").Value = 883.38464
$ws.Range("K116").Value = 883.38464
$ws.Range("M116").Value = 1410.61536

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1862.381
$ws.Range("I3").Value = 883.38464
$ws.Range("K3").Value = 883.38464
$ws.Range("M3").Value = -769.38464

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 55561624
$ws.Range("I86").Value = 83336020
$ws.Range("J86").Value = 12835.667
$ws.Range("K86").Value = 83336020
$ws.Range("L86").Value = 12835.667
$ws.Range("M86").Value = -83334897
$ws.Range("N86").Value = -15081.667

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 55561624
$ws.Range("I89").Value = 83336020
$ws.Range("J89").Value = 12835.667
$ws.Range("K89").Value = 416680100
$ws.Range("L89").Value = 64178.335
$ws.Range("M89").Value = -416674484
$ws.Range("N89").Value = -75410.33499999999

# BSM row 104
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H104").Value = 55789.332
$ws.Range("J104").Value = 55789.332
$ws.Range("L104").Value = 55789.332
$ws.Range("N104").Value = -62777.332

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1564.6216
$ws.Range("I105").Value = 1467.8695
$ws.Range("J105").Value = 1723.5714
$ws.Range("K105").Value = 1467.8695
$ws.Range("L105").Value = 1723.5714
$ws.Range("M105").Value = 279.1305
$ws.Range("N105").Value = -5217.5714

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 12926.134
$ws.Range("I134").Value = 4778.8
$ws.Range("J134").Value = 29220.8
$ws.Range("K134").Value = 14336.4
$ws.Range("L134").Value = 87662.39999999999
$ws.Range("M134").Value = -11801.4
$ws.Range("N134").Value = -92732.39999999999

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2686.5217
$ws.Range("I31").Value = 1974.2858
$ws.Range("J31").Value = 3794.4443
$ws.Range("K31").Value = 1974.2858
$ws.Range("L31").Value = 3794.4443
$ws.Range("M31").Value = -1679.2858
$ws.Range("N31").Value = -4384.4443

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2686.5217
$ws.Range("I34").Value = 1974.2858
$ws.Range("J34").Value = 3794.4443
$ws.Range("K34").Value = 1974.2858
$ws.Range("L34").Value = 3794.4443
$ws.Range("M34").Value = -1772.2858
$ws.Range("N34").Value = -4198.4443

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1820
$ws.Range("I122").Value = 1820
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5460
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3010
$ws.Range("N122").ClearContents()

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 30000.428
$ws.Range("J37").Value = 30000.428
$ws.Range("L37").Value = 90001.284
$ws.Range("N37").Value = -90225.284

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 661.3333
$ws.Range("I98").Value = 180
$ws.Range("J98").Value = 902
$ws.Range("K98").Value = 540
$ws.Range("L98").Value = 2706
$ws.Range("M98").Value = 958
$ws.Range("N98").Value = -5702

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 36.7
$ws.Range("I2").Value = 38.333332
$ws.Range("J2").Value = 22
$ws.Range("K2").Value = 38.333332
$ws.Range("L2").Value = 22
$ws.Range("M2").Value = 74.666668
$ws.Range("N2").Value = -248

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15209687
$ws.Range("I70").Value = 25574482
$ws.Range("J70").Value = 7986.6
$ws.Range("K70").Value = 25574482
$ws.Range("L70").Value = 7986.6
$ws.Range("M70").Value = -25574212
$ws.Range("N70").Value = -8526.6

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 15209687
$ws.Range("I73").Value = 25574482
$ws.Range("J73").Value = 7986.6
$ws.Range("K73").Value = 25574482
$ws.Range("L73").Value = 7986.6
$ws.Range("M73").Value = -25573546
$ws.Range("N73").Value = -9858.6

# GSM row 82
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 30000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 30000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 30000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -30766

# GSM row 85
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 30000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 30000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 30000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -32652

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2127.1724
$ws.Range("I122").Value = 2156.5625
$ws.Range("J122").Value = 2091
$ws.Range("K122").Value = 6469.6875
$ws.Range("L122").Value = 6273
$ws.Range("M122").Value = -4019.6875
$ws.Range("N122").Value = -11173

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1015.7368
$ws.Range("I100").Value = 1064.6364
$ws.Range("J100").Value = 948.5
$ws.Range("K100").Value = 2129.2728
$ws.Range("L100").Value = 1897
$ws.Range("M100").Value = -1588.2728
$ws.Range("N100").Value = -2979

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1064.5927
$ws.Range("I122").Value = 967.2632
$ws.Range("J122").Value = 1295.75
$ws.Range("K122").Value = 2901.7896
$ws.Range("L122").Value = 3887.25
$ws.Range("M122").Value = -451.7896000000001
$ws.Range("N122").Value = -8787.25

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15640699
$ws.Range("I132").Value = 19533514
$ws.Range("J132").Value = 1799580.1
$ws.Range("K132").Value = 58600542
$ws.Range("L132").Value = 5398740.300000001
$ws.Range("M132").Value = -58598012
$ws.Range("N132").Value = -5403800.300000001
